$wb = $excel.ActiveWorkbook

# ---- Update weekly schedule grids for Turma A-E ----
$sheetName = "Turma A"
$ws = $wb.Worksheets.Item($sheetName)
$data = @(
    @("Carlos (Educação Física)", "Pedro (História)", "Beatriz (Inglês)", "Jorge (Geografia)", "Pedro (Português)"),
    @("Alan (Matemática)", "Camila (Artes)", "Bignicius (Ciências)", "Pedro (Português)", "Pedro (História)"),
    @("Bignicius (Ciências)", "Bignicius (Ciências)", "Beatriz (Inglês)", "Alan (Matemática)", "Jorge (Geografia)"),
    @("Alan (Matemática)", "Camila (Artes)", "Bignicius (Ciências)", "Beatriz (Inglês)", "Alan (Matemática)"),
    @("Beatriz (Inglês)", "Pedro (Português)", "Pedro (Português)", "Carlos (Educação Física)", "Carlos (Educação Física)"),
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}
Write-Host ("Updated " + $sheetName)

$sheetName = "Turma B"
$ws = $wb.Worksheets.Item($sheetName)
$data = @(
    @("Beatriz (Inglês)", "Alan (Matemática)", "Jorge (Geografia)", "Bignicius (Ciências)", "Carlos (Educação Física)"),
    @("Alan (Matemática)", "Camila (Artes)", "Alan (Matemática)", "Carlos (Educação Física)", "Alan (Matemática)"),
    @("Pedro (Português)", "Beatriz (Inglês)", "Jorge (Geografia)", "Bignicius (Ciências)", "Carlos (Educação Física)"),
    @("Bignicius (Ciências)", "Bignicius (Ciências)", "Pedro (Português)", "Pedro (Português)", "Camila (Artes)"),
    @("Pedro (Português)", "Pedro (História)", "Beatriz (Inglês)", "Beatriz (Inglês)", "Pedro (História)"),
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}
Write-Host ("Updated " + $sheetName)

$sheetName = "Turma C"
$ws = $wb.Worksheets.Item($sheetName)
$data = @(
    @("Pedro (História)", "Pedro (História)", "Carlos (Educação Física)", "Beatriz (Inglês)", "Beatriz (Inglês)"),
    @("Bignicius (Ciências)", "Alan (Matemática)", "Bignicius (Ciências)", "Alan (Matemática)", "Carlos (Educação Física)"),
    @("Alan (Português)", "Camila (Artes)", "Alan (Português)", "Bignicius (Ciências)", "Alan (Português)"),
    @("Alan (Matemática)", "Alan (Matemática)", "Bignicius (Ciências)", "Alan (Português)", "Camila (Artes)"),
    @("Beatriz (Inglês)", "Carlos (Educação Física)", "Jorge (Geografia)", "Beatriz (Inglês)", "Jorge (Geografia)"),
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}
Write-Host ("Updated " + $sheetName)

$sheetName = "Turma D"
$ws = $wb.Worksheets.Item($sheetName)
$data = @(
    @("Carlos (Educação Física)", "Jorge (Geografia)", "Bignicius (Ciências)", "Pedro (História)", "Beatriz (Inglês)"),
    @("Alan (Matemática)", "Bignicius (Ciências)", "Pedro (História)", "Alan (Português)", "Alan (Português)"),
    @("Jorge (Geografia)", "Alan (Português)", "Alan (Matemática)", "Alan (Matemática)", "Camila (Artes)"),
    @("Alan (Matemática)", "Bignicius (Ciências)", "Camila (Artes)", "Carlos (Educação Física)", "Bignicius (Ciências)"),
    @("Alan (Português)", "Beatriz (Inglês)", "Beatriz (Inglês)", "Beatriz (Inglês)", "Carlos (Educação Física)"),
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}
Write-Host ("Updated " + $sheetName)

$sheetName = "Turma E"
$ws = $wb.Worksheets.Item($sheetName)
$data = @(
    @("Alan (Matemática)", "Carlos (Educação Física)", "Alan (Português)", "Carlos (Educação Física)", "Thiago (Ciências)"),
    @("Beatriz (Inglês)", "Camila (Artes)", "Alan (Matemática)", "Thiago (Ciências)", "Alan (Português)"),
    @("Carlos (Educação Física)", "Alan (Matemática)", "Alan (Português)", "Beatriz (Inglês)", "Pedro (História)"),
    @("Beatriz (Inglês)", "Thiago (Ciências)", "Alan (Matemática)", "Camila (Artes)", "Thiago (Ciências)"),
    @("Pedro (História)", "Beatriz (Inglês)", "Alan (Português)", "Jorge (Geografia)", "Jorge (Geografia)"),
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $row[$c]
    }
}
Write-Host ("Updated " + $sheetName)


# ---- Restructure Workload_Teachers sheet ----
$ws = $wb.Worksheets.Item("Workload_Teachers")

# Style source: B1 currently carries the bold/border/center header style (s=1)
$styleSource = $ws.Cells.Item(1,2)

# New header row: B1:I1 teacher names, J1 = "Ocupação"; old K1 is dropped later
$headers = @("Alan","Bignicius","Jorge","Camila","Thiago","Pedro","Carlos","Beatriz")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $cell = $ws.Cells.Item(1, $col)
    $styleSource.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $headers[$i]
}

$cellJ1 = $ws.Cells.Item(1, 10)
$styleSource.Copy()
$cellJ1.PasteSpecial(-4122)
$cellJ1.Value = "Ocupação"

# Drop the now-unused column K entirely
$ws.Columns.Item(11).Delete()

# Clear the old "Ocupação" data row (numbers used to live in B2:K2)
$ws.Range("B2:K2").ClearContents()

# Style source for column A labels (A2 already carries style s=1)
$occStyleSource = $ws.Cells.Item(2,1)

# New rows: one row per teacher, name in column A, occupation total in column J
$names = @("Alan","Bignicius","Jorge","Camila","Thiago","Pedro","Carlos","Beatriz")
$occ = @(32,16,10,10,4,18,15,20)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = 2 + $i
    $acell = $ws.Cells.Item($r, 1)
    if ($r -ne 2) {
        $occStyleSource.Copy()
        $acell.PasteSpecial(-4122)
    }
    $acell.Value = $names[$i]
    $ws.Cells.Item($r, 10).Value = $occ[$i]
}

Write-Host "Updated Workload_Teachers"
